# Applies the "CDC moratorium elimination" edit described in the commit:
#   - merges the two runs in the "If you make a new agreement..." paragraph
#     into a single run (identical text, just collapsed formatting-run split)
#   - merges the two runs in the "If your landlord agrees..." paragraph
#     into a single run (identical text, just collapsed formatting-run split)
#   - removes the whole CDC-moratorium block of paragraphs (the "You may be
#     eligible..." heading through the "See: " paragraph)
#   - strips the hyperlink/run content out of the paragraph that follows
#     (the one carrying the two-column sectPr + MassLegalHelp.org links),
#     leaving only its paragraph properties, and re-homes the "_GoBack"
#     bookmark (which used to sit in the deleted block) onto that now-empty
#     paragraph.

$d = $word.ActiveDocument

# --- 1. Collapse the "If you make a new agreement..." paragraph's two runs
#        into one by re-"finding & replacing" its full text with itself;
#        Word's Find/Replace naturally re-serializes a single matching run. ---
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute(
    "If you make a new agreement you can have your new agreement replace the old one. This way, it will be enforced the same way as the prior agreement.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If you make a new agreement you can have your new agreement replace the old one. This way, it will be enforced the same way as the prior agreement.",
    2
) | Out-Null

# --- 2. Same trick for the "If your landlord agrees..." paragraph. ---
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Find.Execute(
    "If your landlord agrees to change your agreement, you should still file a Motion to Amend Agreement so that your new agreement can be enforced",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If your landlord agrees to change your agreement, you should still file a Motion to Amend Agreement so that your new agreement can be enforced",
    2
) | Out-Null

# --- 3. Delete the entire CDC-moratorium block: from the "You may be
#        eligible..." Heading3 paragraph through the "See: " paragraph
#        (7 whole paragraphs), inclusive. ---
$blockStart = $d.Paragraphs.Item(9)
$blockEnd = $d.Paragraphs.Item(15)
$blockRange = $d.Range($blockStart.Range.Start, $blockEnd.Range.End)
$blockRange.Delete()

# --- 4. Re-home the "_GoBack" bookmark (it lived inside the block we just
#        deleted) onto the start of the paragraph that used to hold the
#        MassLegalHelp.org hyperlinks -- must happen BEFORE we clear that
#        paragraph's content below. ---
$linkPara = $d.Paragraphs.Item(9)
$bookmarkRange = $d.Range($linkPara.Range.Start, $linkPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# --- 5. Strip the hyperlinks/runs out of that paragraph, keeping its
#        paragraph mark (and thus its pPr/sectPr) intact. ---
$linkPara = $d.Paragraphs.Item(9)
$contentRange = $d.Range($linkPara.Range.Start, $linkPara.Range.End - 1)
$contentRange.Delete()
